$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in the existing "RequireItem/plant quest" description (H11):
# "G|奥莱伊李交了你种植植物的方法..." -> "G|奥莱伊李||教了你种植植物的方法..."
$ws.Range("H11").Value = "G|奥莱伊李||教了你种植植物的方法，||你需要到附近的田地里，种下|Y|豌豆种子||。并在收获后，把果实带交给他。"

# Add a brand new quest row (row 12: Id 12000009) - "教训猩猩" (Teach the ape a lesson)
$ws.Range("B12").Value = "教训猩猩"
$ws.Range("C12").Value = "onelesson"
$ws.Range("D12").Value = 12000005
$ws.Range("H12").Value = "G|塞巴斯恰恩||委托你找到可恶的|G|科迪||，并通过卡牌战斗的方式击败他。|G|科迪|经常欺负周边的村民，非常可恶。"
$ws.Range("I12").Value = 42120003

# J12 (Type = "增项") reuses the same highlight formatting already used for
# this quest-type value elsewhere in the column (e.g. J5), so copy that
# cell's format before writing the value.
$ws.Range("J5").Copy()
$ws.Range("J12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J12").Value = "增项"

$ws.Range("M12").Value = "npckedi"
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = "questonelesson"
$ws.Range("S12").Value = 50
$ws.Range("AA12").Value = 43020103

# Adjust the display-x value of the following quest row (F13)
$ws.Range("F13").Value = 15
